# moar data - wjec 2020
# Add a new ingest-log row (row 8) for the WJEC 2020 exam timetable PDF:
#   - C8 gets the source URL (as a hyperlink, styled like the other source cells)
#   - D8 gets the "Accessed" date (dd/mm/yyyy;@), same as the other rows
# Also moves the saved cursor selection to C8 (where the new data was entered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.wjec.co.uk/exam-officers/Examination_Timetable_2020_Final_update_041219%20(2).pdf?language_id=1"

# --- D8: accessed date, matching the date format/style used by the rest of column D ---
$ws.Range("D8").Value = 43863
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)   # xlPasteFormats - copy over just the number format/style

# --- C8: the new source link, registered as a real hyperlink ---
$ws.Range("C8").Value = $url
$null = $ws.Hyperlinks.Add($ws.Range("C8"), $url)

# Hyperlinks.Add stamps its own (fresh) style on the cell; restore the shared
# "Hyperlink" look used by the other source cells (C3/C6/C7) instead of leaving
# a brand-new, duplicate style behind.
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Move the active selection to the newly-entered cell.
$null = $ws.Range("C8").Select()
